$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1921.875
$ws.Range("J17").Value = 1921.875
$ws.Range("L17").Value = 5765.625
$ws.Range("N17").Value = -6101.625
$ws.Range("H32").Value = 2546.5881
$ws.Range("I32").Value = 1723
$ws.Range("K32").Value = 1723
$ws.Range("M32").Value = -1397
$ws.Range("H40").Value = 4996.5835
$ws.Range("I40").Value = 3665.3333
$ws.Range("J40").Value = 5440.3335
$ws.Range("K40").Value = 3665.3333
$ws.Range("L40").Value = 5440.3335
$ws.Range("M40").Value = -3490.3333
$ws.Range("N40").Value = -5790.3335
$ws.Range("H74").Value = 5199
$ws.Range("I74").Value = 4358.7
$ws.Range("J74").Value = 8000
$ws.Range("K74").Value = 4358.7
$ws.Range("L74").Value = 8000
$ws.Range("M74").Value = -3422.7
$ws.Range("N74").Value = -9872
$ws.Range("H77").Value = 5199
$ws.Range("I77").Value = 4358.7
$ws.Range("J77").Value = 8000
$ws.Range("K77").Value = 21793.5
$ws.Range("L77").Value = 40000
$ws.Range("M77").Value = -17113.5
$ws.Range("N77").Value = -49360
$ws.Range("H97").Value = 1697.2
$ws.Range("J97").Value = 1946.5
$ws.Range("L97").Value = 5839.5
$ws.Range("N97").Value = -6831.5
$ws.Range("H137").Value = 4233
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 4233
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 12699
$ws.Range("M137").ClearContents()
$ws.Range("N137").Value = -17799
$ws.Range("H138").Value = 1444.9474
$ws.Range("I138").Value = 1444.9474
$ws.Range("K138").Value = 4334.8422
$ws.Range("M138").Value = 805.1578

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 842.375
$ws.Range("I2").Value = 677
$ws.Range("J2").Value = 2000
$ws.Range("K2").Value = 677
$ws.Range("L2").Value = 2000
$ws.Range("M2").Value = -564
$ws.Range("N2").Value = -2226
$ws.Range("H5").Value = 128.83333
$ws.Range("I5").Value = 140.8
$ws.Range("K5").Value = 140.8
$ws.Range("M5").Value = -28.80000000000001
$ws.Range("H45").Value = 1846.25
$ws.Range("I45").Value = 1996.8
$ws.Range("K45").Value = 1996.8
$ws.Range("M45").Value = -1619.8
$ws.Range("H74").Value = 1531.75
$ws.Range("I74").Value = 1428.1
$ws.Range("J74").Value = 2050
$ws.Range("K74").Value = 1428.1
$ws.Range("L74").Value = 2050
$ws.Range("M74").Value = -554.0999999999999
$ws.Range("N74").Value = -3798
$ws.Range("H77").Value = 1531.75
$ws.Range("I77").Value = 1428.1
$ws.Range("J77").Value = 2050
$ws.Range("K77").Value = 7140.5
$ws.Range("L77").Value = 10250
$ws.Range("M77").Value = -2772.5
$ws.Range("N77").Value = -18986
$ws.Range("H97").Value = 1266.5883
$ws.Range("I97").Value = 720.75
$ws.Range("K97").Value = 720.75
$ws.Range("M97").Value = -224.75
$ws.Range("H116").Value = 842.375
$ws.Range("I116").Value = 677
$ws.Range("J116").Value = 2000
$ws.Range("K116").Value = 677
$ws.Range("L116").Value = 2000
$ws.Range("M116").Value = 1617
$ws.Range("N116").Value = -6588

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 842.375
$ws.Range("I3").Value = 677
$ws.Range("J3").Value = 2000
$ws.Range("K3").Value = 677
$ws.Range("L3").Value = 2000
$ws.Range("M3").Value = -563
$ws.Range("N3").Value = -2228
$ws.Range("H4").Value = 128.83333
$ws.Range("I4").Value = 140.8
$ws.Range("K4").Value = 140.8
$ws.Range("M4").Value = -25.80000000000001
$ws.Range("H107").Value = 2242
$ws.Range("I107").Value = 1485.5
$ws.Range("K107").Value = 1485.5
$ws.Range("M107").Value = 434.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2046.9231
$ws.Range("I31").Value = 1833.1111
$ws.Range("J31").Value = 2528
$ws.Range("K31").Value = 1833.1111
$ws.Range("L31").Value = 2528
$ws.Range("M31").Value = -1538.1111
$ws.Range("N31").Value = -3118
$ws.Range("H34").Value = 2046.9231
$ws.Range("I34").Value = 1833.1111
$ws.Range("J34").Value = 2528
$ws.Range("K34").Value = 1833.1111
$ws.Range("L34").Value = 2528
$ws.Range("M34").Value = -1631.1111
$ws.Range("N34").Value = -2932
$ws.Range("H58").Value = 4892.4443
$ws.Range("I58").Value = 2895
$ws.Range("J58").Value = 6490.4
$ws.Range("K58").Value = 2895
$ws.Range("L58").Value = 6490.4
$ws.Range("M58").Value = -2692
$ws.Range("N58").Value = -6896.4
$ws.Range("H99").Value = 3789.6667
$ws.Range("I99").Value = 3729.7144
$ws.Range("J99").Value = 3999.5
$ws.Range("K99").Value = 3729.7144
$ws.Range("L99").Value = 3999.5
$ws.Range("M99").Value = -2231.7144
$ws.Range("N99").Value = -6995.5
$ws.Range("H126").Value = 3789.6667
$ws.Range("I126").Value = 3729.7144
$ws.Range("J126").Value = 3999.5
$ws.Range("K126").Value = 11189.1432
$ws.Range("L126").Value = 11998.5
$ws.Range("M126").Value = -8719.143199999999
$ws.Range("N126").Value = -16938.5
$ws.Range("H134").Value = 3685.5881
$ws.Range("I134").Value = 3650.3333
$ws.Range("J134").Value = 3950
$ws.Range("K134").Value = 10950.9999
$ws.Range("L134").Value = 11850
$ws.Range("M134").Value = -8415.999899999999
$ws.Range("N134").Value = -16920
$ws.Range("H136").Value = 4892.4443
$ws.Range("I136").Value = 2895
$ws.Range("J136").Value = 6490.4
$ws.Range("K136").Value = 8685
$ws.Range("L136").Value = 19471.2
$ws.Range("M136").Value = -6135
$ws.Range("N136").Value = -24571.2

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 5181.125
$ws.Range("J69").Value = 5600
$ws.Range("L69").Value = 16800
$ws.Range("N69").Value = -18422
$ws.Range("H70").Value = 12935.5
$ws.Range("I70").Value = 1161.3334
$ws.Range("K70").Value = 3484.0002
$ws.Range("M70").Value = -3169.0002
$ws.Range("H72").Value = 5181.125
$ws.Range("J72").Value = 5600
$ws.Range("L72").Value = 50400
$ws.Range("N72").Value = -58512
$ws.Range("H73").Value = 12935.5
$ws.Range("I73").Value = 1161.3334
$ws.Range("K73").Value = 3484.0002
$ws.Range("M73").Value = -2392.0002
$ws.Range("H75").Value = 732.1667
$ws.Range("H78").Value = 732.1667
$ws.Range("H132").Value = 1723.875
$ws.Range("I132").Value = 1730.6666
$ws.Range("J132").Value = 1719.8
$ws.Range("K132").Value = 15575.9994
$ws.Range("L132").Value = 15478.2
$ws.Range("M132").Value = -13045.9994
$ws.Range("N132").Value = -20538.2

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4011.6667
$ws.Range("I102").Value = 4012.818
$ws.Range("K102").Value = 4012.818
$ws.Range("M102").Value = -2390.818
$ws.Range("H132").Value = 3400
$ws.Range("I132").Value = 3400
$ws.Range("K132").Value = 10200
$ws.Range("M132").Value = -7670

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3547.9565
$ws.Range("I7").Value = 3200.1904
$ws.Range("K7").Value = 3200.1904
$ws.Range("M7").Value = -3088.1904
$ws.Range("H40").Value = 1396.5714
$ws.Range("I40").Value = 1396.5714
$ws.Range("K40").Value = 1396.5714
$ws.Range("M40").Value = -1260.5714
$ws.Range("H46").Value = 1266.4546
$ws.Range("J46").Value = 1349
$ws.Range("L46").Value = 1349
$ws.Range("N46").Value = -1725
$ws.Range("H126").Value = 3547.9565
$ws.Range("I126").Value = 3200.1904
$ws.Range("K126").Value = 9600.5712
$ws.Range("M126").Value = -7130.5712
$ws.Range("H136").Value = 5922.5835
$ws.Range("I136").Value = 3257.1
$ws.Range("K136").Value = 9771.299999999999
$ws.Range("M136").Value = -7221.299999999999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1258.8823
$ws.Range("I126").Value = 1224.9231
$ws.Range("J126").Value = 1369.25
$ws.Range("K126").Value = 3674.7693
$ws.Range("L126").Value = 4107.75
$ws.Range("M126").Value = -1204.7693
$ws.Range("N126").Value = -9047.75
